$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()
$ws.Rows.Item(49).Select()
$ws.Rows.Item(49).Delete()
